# Femacal de La Calera - Berenjena: weekly price update.
# A new weekly record is inserted above the existing row 431, pushing rows
# 431-486 down to 432-487 (dimension grows from A1:R486 to A1:R487).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 431 (shifts 431:486 -> 432:487).
$ws.Rows.Item(431).Insert()

# Populate the newly inserted row 431 with the latest record.
$ws.Cells.Item(431, 1).Value = 3
$ws.Cells.Item(431, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(431, 3).Value = "Coquimbo"
$ws.Cells.Item(431, 4).Value = 45124
$ws.Cells.Item(431, 5).Value = 5
$ws.Cells.Item(431, 6).Value = 100112001
$ws.Cells.Item(431, 7).Value = "Berenjena"
$ws.Cells.Item(431, 8).Value = "Sin especificar"
$ws.Cells.Item(431, 9).Value = "Primera"
$ws.Cells.Item(431, 10).Value = 100
$ws.Cells.Item(431, 11).Value = 7000
$ws.Cells.Item(431, 12).Value = 7500
$ws.Cells.Item(431, 13).Value = 7250
$ws.Cells.Item(431, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(431, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(431, 16).Value = 121
$ws.Cells.Item(431, 17).Value = 60
$ws.Cells.Item(431, 18).Value = "Hortaliza"
